$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Use a scratch cell far outside the used range to produce a text value
# ("11.02.2021") that Excel will not auto-convert into a date serial
# number. A formula that evaluates to the literal string, copied and
# pasted as values into the destination cells, lands as a plain shared
# string without picking up (or creating) any extra number-format style.
$scratch = $ws.Cells.Item(500, 500)
$scratch.Formula = "=""11.02.2021"""
$scratch.Copy()

$ws.Cells.Item(3, 1).PasteSpecial(-4163)
$ws.Cells.Item(4, 1).PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = $false

$ws.Cells.Item(3, 2).Value = "Morgan"
$ws.Cells.Item(3, 3).Value = "Pressed R,D in quick succession and snake went up, then no controls worked until the game was over"

$ws.Cells.Item(4, 2).Value = "Morgan"
$ws.Cells.Item(4, 3).Value = "Due to the fact that the gridding is off for food and snake navigation, snake would eat food that was at least 1 snake width away from the snake. Need to work on accuracy and making the game a consistent square grid?"

$ws.Range("A5").Select()
